$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 1.186041333333334
$ws.Range("H2").Value = 3.558124
$ws.Range("I2").Value = 0.3004487152423012
$ws.Range("J2").Value = 0.3004487152423012
$ws.Range("M2").Value = 0.4349816666666667
$ws.Range("N2").Value = 1.304945
$ws.Range("O2").Value = 0.0068619340787224
$ws.Range("P2").Value = 0.00687614195861717
$ws.Range("Q2").Value = 0.515906235908889
$ws.Range("R2").Value = 4.643156123180001
$ws.Range("S2").Value = 0.002061659278029509
$ws.Range("T2").Value = 0.002065928017290209
$ws.Range("G3").Value = 1.186041333333334
$ws.Range("H3").Value = 3.558124
$ws.Range("I3").Value = 0.3004487152423012
$ws.Range("J3").Value = 0.3004487152423012
$ws.Range("O3").Value = 0.005722841821244052
$ws.Range("P3").Value = 0.00573469116988544
$ws.Range("Q3").Value = 0.4302649586586667
$ws.Range("R3").Value = 3.872384627928001
$ws.Range("S3").Value = 0.001719420472727687
$ws.Range("T3").Value = 0.00172298059430345
$ws.Range("G4").Value = 1.186041333333334
$ws.Range("H4").Value = 3.558124
$ws.Range("I4").Value = 0.3004487152423012
$ws.Range("J4").Value = 0.3004487152423012
$ws.Range("M4").Value = 29.338587
$ws.Range("N4").Value = 88.015761
$ws.Range("O4").Value = 0.4628228391775791
$ws.Range("P4").Value = 0.4637811304167767
$ws.Range("Q4").Value = 34.79677684359601
$ws.Range("R4").Value = 313.170991592364
$ws.Range("S4").Value = 0.1390545274156978
$ws.Range("T4").Value = 0.1393424447873427
$ws.Range("G5").Value = 1.186041333333334
$ws.Range("H5").Value = 3.558124
$ws.Range("I5").Value = 0.3004487152423012
$ws.Range("J5").Value = 0.3004487152423012
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.3929435
$ws.Range("N5").Value = 0.785887
$ws.Range("O5").Value = 0.006198772500747056
$ws.Range("P5").Value = 0.004141071520586516
$ws.Range("Q5").Value = 0.4660472326646667
$ws.Range("R5").Value = 2.796283395988
$ws.Range("S5").Value = 0.001862413233928759
$ws.Range("T5").Value = 0.001244179618086702
$ws.Range("G6").Value = 1.186041333333334
$ws.Range("H6").Value = 3.558124
$ws.Range("I6").Value = 0.3004487152423012
$ws.Range("J6").Value = 0.3004487152423012
$ws.Range("M6").Value = 32.861248
$ws.Range("N6").Value = 98.583744
$ws.Range("O6").Value = 0.5183936124217073
$ws.Range("P6").Value = 0.5194669649341341
$ws.Range("Q6").Value = 38.97479839291734
$ws.Range("R6").Value = 350.773185536256
$ws.Range("S6").Value = 0.1557506948419174
$ws.Range("T6").Value = 0.1560731822252781
$ws.Range("I7").Value = 0.331325035675986
$ws.Range("J7").Value = 0.3313250356759861
$ws.Range("M7").Value = 0.4349816666666667
$ws.Range("N7").Value = 1.304945
$ws.Range("O7").Value = 0.0068619340787224
$ws.Range("P7").Value = 0.00687614195861717
$ws.Range("Q7").Value = 0.5689245563261111
$ws.Range("R7").Value = 5.120321006935
$ws.Range("S7").Value = 0.002273530553438964
$ws.Range("T7").Value = 0.002278237979751979
$ws.Range("I8").Value = 0.331325035675986
$ws.Range("J8").Value = 0.3313250356759861
$ws.Range("O8").Value = 0.005722841821244052
$ws.Range("P8").Value = 0.00573469116988544
$ws.Range("S8").Value = 0.001896120770591711
$ws.Range("T8").Value = 0.001900046756453056
$ws.Range("I9").Value = 0.331325035675986
$ws.Range("J9").Value = 0.3313250356759861
$ws.Range("M9").Value = 29.338587
$ws.Range("N9").Value = 88.015761
$ws.Range("O9").Value = 0.4628228391775791
$ws.Range("P9").Value = 0.4637811304167767
$ws.Range("Q9").Value = 38.372749638207
$ws.Range("R9").Value = 345.354746743863
$ws.Range("S9").Value = 0.1533447937021725
$ws.Range("T9").Value = 0.1536622995811877
$ws.Range("I10").Value = 0.331325035675986
$ws.Range("J10").Value = 0.3313250356759861
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.3929435
$ws.Range("N10").Value = 0.785887
$ws.Range("O10").Value = 0.006198772500747056
$ws.Range("P10").Value = 0.004141071520586516
$ws.Range("Q10").Value = 0.5139416750868333
$ws.Range("R10").Value = 3.083650050521
$ws.Range("S10").Value = 0.00205380851995734
$ws.Range("T10").Value = 0.001372040669295137
$ws.Range("I11").Value = 0.331325035675986
$ws.Range("J11").Value = 0.3313250356759861
$ws.Range("M11").Value = 32.861248
$ws.Range("N11").Value = 98.583744
$ws.Range("O11").Value = 0.5183936124217073
$ws.Range("P11").Value = 0.5194669649341341
$ws.Range("Q11").Value = 42.98013542039466
$ws.Range("R11").Value = 386.821218783552
$ws.Range("S11").Value = 0.1717567821298255
$ws.Range("T11").Value = 0.1721124106892982
$ws.Range("G12").Value = 0.6500023333333333
$ws.Range("H12").Value = 1.950007
$ws.Range("I12").Value = 0.1646589882374797
$ws.Range("J12").Value = 0.1646589882374797
$ws.Range("M12").Value = 0.4349816666666667
$ws.Range("N12").Value = 1.304945
$ws.Range("O12").Value = 0.0068619340787224
$ws.Range("P12").Value = 0.00687614195861717
$ws.Range("Q12").Value = 0.2827390982905555
$ws.Range("R12").Value = 2.544651884615
$ws.Range("S12").Value = 0.001129879122754713
$ws.Range("T12").Value = 0.001132218577883185
$ws.Range("G13").Value = 0.6500023333333333
$ws.Range("H13").Value = 1.950007
$ws.Range("I13").Value = 0.1646589882374797
$ws.Range("J13").Value = 0.1646589882374797
$ws.Range("O13").Value = 0.005722841821244052
$ws.Range("P13").Value = 0.00573469116988544
$ws.Range("Q13").Value = 0.2358039464726667
$ws.Range("R13").Value = 2.122235518254
$ws.Range("S13").Value = 0.000942317344129181
$ws.Range("T13").Value = 0.0009442684458877451
$ws.Range("G14").Value = 0.6500023333333333
$ws.Range("H14").Value = 1.950007
$ws.Range("I14").Value = 0.1646589882374797
$ws.Range("J14").Value = 0.1646589882374797
$ws.Range("M14").Value = 29.338587
$ws.Range("N14").Value = 88.015761
$ws.Range("O14").Value = 0.4628228391775791
$ws.Range("P14").Value = 0.4637811304167767
$ws.Range("Q14").Value = 19.070150006703
$ws.Range("R14").Value = 171.631350060327
$ws.Range("S14").Value = 0.07620794043217793
$ws.Range("T14").Value = 0.07636573169806106
$ws.Range("G15").Value = 0.6500023333333333
$ws.Range("H15").Value = 1.950007
$ws.Range("I15").Value = 0.1646589882374797
$ws.Range("J15").Value = 0.1646589882374797
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.5
$ws.Range("M15").Value = 0.3929435
$ws.Range("N15").Value = 0.785887
$ws.Range("O15").Value = 0.006198772500747056
$ws.Range("P15").Value = 0.004141071520586516
$ws.Range("Q15").Value = 0.2554141918681667
$ws.Range("R15").Value = 1.532485151209
$ws.Range("S15").Value = 0.001020683608287322
$ws.Range("T15").Value = 0.0006818646467988172
$ws.Range("G16").Value = 0.6500023333333333
$ws.Range("H16").Value = 1.950007
$ws.Range("I16").Value = 0.1646589882374797
$ws.Range("J16").Value = 0.1646589882374797
$ws.Range("M16").Value = 32.861248
$ws.Range("N16").Value = 98.583744
$ws.Range("O16").Value = 0.5183936124217073
$ws.Range("P16").Value = 0.5194669649341341
$ws.Range("Q16").Value = 21.35988787624533
$ws.Range("R16").Value = 192.238990886208
$ws.Range("S16").Value = 0.0853581677301305
$ws.Range("T16").Value = 0.08553490486884885
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8035953333333333
$ws.Range("H17").Value = 2.410786
$ws.Range("I17").Value = 0.2035672608442332
$ws.Range("J17").Value = 0.2035672608442332
$ws.Range("M17").Value = 0.4349816666666667
$ws.Range("N17").Value = 1.304945
$ws.Range("O17").Value = 0.0068619340787224
$ws.Range("P17").Value = 0.00687614195861717
$ws.Range("Q17").Value = 0.3495492374188889
$ws.Range("R17").Value = 3.14594313677
$ws.Range("S17").Value = 0.001396865124499216
$ws.Range("T17").Value = 0.001399757383691798
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.8035953333333333
$ws.Range("H18").Value = 2.410786
$ws.Range("I18").Value = 0.2035672608442332
$ws.Range("J18").Value = 0.2035672608442332
$ws.Range("O18").Value = 0.005722841821244052
$ws.Range("P18").Value = 0.00573469116988544
$ws.Range("Q18").Value = 0.2915234934546667
$ws.Range("R18").Value = 2.623711441092
$ws.Range("S18").Value = 0.001164983233795474
$ws.Range("T18").Value = 0.00116739537324119
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.8035953333333333
$ws.Range("H19").Value = 2.410786
$ws.Range("I19").Value = 0.2035672608442332
$ws.Range("J19").Value = 0.2035672608442332
$ws.Range("M19").Value = 29.338587
$ws.Range("N19").Value = 88.015761
$ws.Range("O19").Value = 0.4628228391775791
$ws.Range("P19").Value = 0.4637811304167767
$ws.Range("Q19").Value = 23.576351599794
$ws.Range("R19").Value = 212.187164398146
$ws.Range("S19").Value = 0.09421557762753081
$ws.Range("T19").Value = 0.09441065435018529
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.8035953333333333
$ws.Range("H20").Value = 2.410786
$ws.Range("I20").Value = 0.2035672608442332
$ws.Range("J20").Value = 0.2035672608442332
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.5
$ws.Range("M20").Value = 0.3929435
$ws.Range("N20").Value = 0.785887
$ws.Range("O20").Value = 0.006198772500747056
$ws.Range("P20").Value = 0.004141071520586516
$ws.Range("Q20").Value = 0.3157675628636666
$ws.Range("R20").Value = 1.894605377182
$ws.Range("S20").Value = 0.001261867138573635
$ws.Range("T20").Value = 0.0008429865864058606
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.8035953333333333
$ws.Range("H21").Value = 2.410786
$ws.Range("I21").Value = 0.2035672608442332
$ws.Range("J21").Value = 0.2035672608442332
$ws.Range("M21").Value = 32.861248
$ws.Range("N21").Value = 98.583744
$ws.Range("O21").Value = 0.5183936124217073
$ws.Range("P21").Value = 0.5194669649341341
$ws.Range("Q21").Value = 26.40714554030933
$ws.Range("R21").Value = 237.664309862784
$ws.Range("S21").Value = 0.105527967719834
$ws.Range("T21").Value = 0.105746467150709
